$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "Art"

# Row 2: change A2 from Breitschwert to Dolch, add C2 = Waffe
# (write C2 first so shared-string order matches: Art, Waffe, Dolch, ...)
$ws.Range("C2").Value = "Waffe"
$ws.Range("A2").Value = "Dolch"

# Row 3: new row - Lederharnisch / 1 / Rüstung
$ws.Range("A3").Value = "Lederharnisch"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Rüstung"

# Column C width (target stored width 17.88671875; engine snaps to 1/6
# character increments, so 17 is the closest input that lands on the
# nearest achievable stored value, 17.833333...)
$ws.Columns.Item(3).ColumnWidth = 17

# Selection on C3 as per diff
$ws.Range("C3").Select()
